$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.658.05"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").Value = "1.637.47"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.500"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.40%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +2.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0624"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("E11").Value = "  +3.19%  "
$ws.Range("D12").Value = "1.864.33"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").Value = "1.622.65"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("E14").Value = "  +1.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.529"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").Value = "26.669.31"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.39%  "
$ws.Range("D18").Value = "0.0₃0745"
$ws.Range("E18").Value = "  +1.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.69%  "
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("E22").Value = "  +1.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.66%  "
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.03%  "
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("E28").Value = "  +4.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.51%  "
$ws.Range("E30").Value = "  -1.71%  "
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.30"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.41%  "
$ws.Range("E33").Value = "  -0.60%  "
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("E35").Value = "  -2.58%  "
$ws.Range("D36").Value = "1.191.95"
$ws.Range("E36").Value = "  +1.02%  "
$ws.Range("E37").Value = "  +5.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.810"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.508"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.77%  "
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.795"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").Value = "1.772.49"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("E46").Value = "  +1.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.40%  "
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.67"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.08%  "
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("E51").Value = "  +0.11%  "
